# Add two new "hydrogen-fuelled gas power" sector reference rows.
#
# The sector reference table lists each sector code twice (columns B and C)
# with a constant "(new)" marker in column D. Two new sector codes are being
# introduced, each inserted directly after its corresponding "*_gas_power_ccs"
# row:
#   - "18_01_02_gas_power_h2" goes right after "18_01_02_gas_power_ccs"
#   - "09_01_02_gas_power_h2" goes right after "09_01_02_gas_power_ccs"
#
# The "18_..." row sits further down the sheet (originally row 249, just
# before "18_01_03_oil"), so it is inserted first so its row number isn't
# disturbed by the earlier ("09_...") insertion that happens afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert "18_01_02_gas_power_h2" after the "18_01_02_gas_power_ccs" row
# (row 248 -> new row pushed in at 249, shifting the old row 249
# "18_01_03_oil" down to row 250).
$ws.Range("A249").EntireRow.Insert()
$ws.Range("B249").Value = "18_01_02_gas_power_h2"
$ws.Range("C249").Value = "18_01_02_gas_power_h2"
$ws.Range("D249").Value = "(new)"

# Insert "09_01_02_gas_power_h2" after the "09_01_02_gas_power_ccs" row
# (row 116 -> new row pushed in at 117, shifting the old row 117
# "09_01_02_01_gasturbine" down to row 118, and everything below it,
# including the row just inserted above, down by one more row).
$ws.Range("A117").EntireRow.Insert()
$ws.Range("B117").Value = "09_01_02_gas_power_h2"
$ws.Range("C117").Value = "09_01_02_gas_power_h2"
$ws.Range("D117").Value = "(new)"

# Match the author's final selection/viewport in the sheet.
$ws.Range("C118").Select() | Out-Null
